# Revert "new changes in ops (ordercreation & orderpage & order form)"
# This reverts the prior edit: restore old order data, old column widths,
# the old "active cell" selection, and the orphaned conditional-format dxf.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: E1/F1 ("Typist" / "Typist QC") keep the exact same text before and
# after this revert - only their underlying shared-string index shifts
# because of the string-table churn below, so no write is needed here.

# --- Row 2 ---
$ws.Range("A2").Value = 45437.0416666088
$ws.Range("B2").Value = "ST18-002"
$ws.Range("J2").Value = "Full Search"

# --- Row 4 ---
$ws.Range("B4").Value = "ST18-008"
$ws.Range("C4").Value = "SIPL5316"
$ws.Range("D4").Value = "SIPL5688"
# E4 / F4 no longer hold any data in the reverted version - remove the
# cells outright (not just clear their contents) so they disappear from
# the sheet entirely, same as the original template.
$ws.Range("E4").Clear() | Out-Null
$ws.Range("F4").Clear() | Out-Null
$ws.Range("M4").Value = "WIP"

# --- Column widths (revert to the narrower, bestFit pre-edit sizes) ---
$ws.Columns("C").ColumnWidth = 35.5
$ws.Columns("G").ColumnWidth = 26.666666666666668
$ws.Columns("J").ColumnWidth = 31.333333333333332

# --- Restore old selection/active cell ---
$ws.Range("A16").Select() | Out-Null

# --- Recreate the orphaned conditional-formatting dxf (light red fill /
#     dark red text) that is left behind in styles.xml once the rule
#     that used to reference it was removed. ---
$rngCf = $ws.Range("B2:B4")
$fc = $rngCf.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = -16383844
$fc.Interior.Color = 13551615
$fc.Delete()
